$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block in columns J/K (rows 3-11) -----------------------------
# Enter the label text in the specific order needed so the shared-string
# table ends up with the same unique-string order as the target workbook.
$ws.Cells.Item(3, 10).Value  = "Desity"
$ws.Cells.Item(5, 10).Value  = "Poisson ratio"
$ws.Cells.Item(11, 10).Value = "Interlaminar shearstrength (MPa)"
$ws.Cells.Item(10, 10).Value = "In-plane shear strength(MPa)"
$ws.Cells.Item(4, 10).Value  = "Young's modulus (GPa)"
$ws.Cells.Item(6, 10).Value  = "In-plane shear modulus(GPa)"
$ws.Cells.Item(7, 10).Value  = "Interlaminar shear modulus(GPa)"
$ws.Cells.Item(8, 10).Value  = "Tensile Strength(MPa)"
$ws.Cells.Item(9, 10).Value  = "Compressive strength(MPa)"

# Known numeric values for density / Young's modulus; the rest are left
# blank (placeholders for values still to be filled in).
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(4, 11).Value = 10.1

# Column J needs to be wide enough to show the longest label.
$ws.Columns.Item(10).ColumnWidth = 31.140625

# --- Thin box border around the whole J3:K11 block --------------------------
$rng = $ws.Range("J3:K11")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# --- Selection now on the new block ----------------------------------------
$rng.Select() | Out-Null

# --- Reposition / resize the picture ----------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 183.75007874015748
$shp.Top = 24
$shp.Width = 246
$shp.Height = 227.07692913385827
